$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the summary header fields -------------------------------------
$ws.Range("E11").Value = 971112   # VALOR MORA
$ws.Range("C13").Value = 2        # Cant. Trabajadores
$ws.Range("F13").Value = 19       # Cant. Periodos

# --- Rebuild the worker detail table ---------------------------------------
# The old table had 20 data rows (16-35): one YESSICA record, one LINDA ROSA
# record, 17 ERIKA records and a final DANIELA record. The new table keeps
# the YESSICA record as-is, drops the LINDA ROSA record and now lists
# ERIKA's account for every period from 2403 through 2412 and 2501 through
# 2508 (18 rows), for a total of 19 data rows (16-34).
#
# Deleting row 17 shifts every following row up by one; because row 35 (the
# old last row, with its own "closing" border style) lands on row 34, the
# special bottom-border formatting of the last table row is preserved
# automatically.
$ws.Rows.Item(17).Delete()

$periods = @("2403","2404","2405","2406","2407","2408","2409","2410","2411","2412","2501","2502","2503","2504","2505","2506","2507","2508")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 17 + $i
    $ws.Cells.Item($r, 3).Value = "1075544728"
    $ws.Cells.Item($r, 4).Value = "ERIKA DANIELA DUSSAN GARCIA"
    $ws.Cells.Item($r, 5).Value = $periods[$i]
    $ws.Cells.Item($r, 6).Value = 52000
    $ws.Cells.Item($r, 7).Value = 1300000
}

# Column D ("Nombre Trabajador") can now be narrower since the longest
# remaining name is shorter than "DANIELA CAROLINA CASTRO CABARCAS" was.
$ws.Columns.Item(4).ColumnWidth = 28.9
